# Daily attendance processing - 2025-12-06 16:29:50
# Normalize the "Recorded By" (column G) entries so that the two
# contributors are listed in their corrected order for the affected
# session rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of exact current "Recorded By" text -> corrected text.
# Only rows whose value matches one of these (case-sensitive) are touched.
$swapMap = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = $cell.Value()

    if ($swapMap.ContainsKey($current)) {
        $cell.Value = $swapMap[$current]
    }
}
